# edit.ps1 - apply the reference-docx formatting/content tweaks described
# by the commit "Edit ms; add plos template; update docx ref doc".
#
# Two kinds of changes happen here:
#   1. Body text: the stray "_GoBack" bookmark that Word leaves behind at
#      the cursor's last-edit position is relocated from the middle of
#      "A paragraph of text with some " (splitting "s" / "ome ") down into
#      the first "And some more text." sentence (splitting "And some" /
#      " more text."), and the two runs it used to separate are re-merged
#      into a single run.
#   2. Styles: Normal/Heading1/Heading2/Heading3 pick up new spacing,
#      justification and font-size tweaks.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1a. Merge "A paragraph of text with s" + bookmark + "ome " back into a
#     single run and drop the old bookmark that split them.
# ---------------------------------------------------------------------
$d.Content.Find.Execute("A paragraph of text with some ", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "A paragraph of text with some ", 2) | Out-Null

# ---------------------------------------------------------------------
# 1b. Re-create the "_GoBack" bookmark, this time splitting the first
#     "And some more text." sentence (the one right after "code.") into
#     "And some" / " more text.".
# ---------------------------------------------------------------------
$text = $d.Content.Text
$idx = $text.IndexOf("And some more text.")
$pos = $idx + ("And some").Length
$r = $d.Range($pos, $pos)
$d.Bookmarks.Add("_GoBack", $r) | Out-Null

# ---------------------------------------------------------------------
# 2. Style tweaks
# ---------------------------------------------------------------------

# Normal: tighter exact line spacing + full justification
$normal = $d.Styles("Normal")
$normal.ParagraphFormat.LineSpacingRule = 4   # wdLineSpaceExactly
$normal.ParagraphFormat.LineSpacing = 22      # 440 twips / 20 = 22 pt
$normal.ParagraphFormat.Alignment = 3         # wdAlignParagraphJustify ("both")

# Heading 1: less space before, slightly smaller
$h1 = $d.Styles("Heading 1")
$h1.ParagraphFormat.SpaceBefore = 6
$h1.Font.Size = 17

# Heading 2: no space before, bold instead of italic, smaller
$h2 = $d.Styles("Heading 2")
$h2.ParagraphFormat.SpaceBefore = 0
$h2.Font.Bold = $true
$h2.Font.Italic = $false
$h2.Font.Size = 13

# Heading 3: less space before
$h3 = $d.Styles("Heading 3")
$h3.ParagraphFormat.SpaceBefore = 6

Write-Host "edit applied"
